$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.049949195306054
$ws.Cells.Item(2, 4).Value = 1.054878083059933
$ws.Cells.Item(2, 5).Value = 1.062921069477074
$ws.Cells.Item(2, 6).Value = 1.068902345652393
$ws.Cells.Item(2, 9).Value = 1.044522586718311
$ws.Cells.Item(2, 10).Value = 1.054984936093648
$ws.Cells.Item(2, 11).Value = 1.057620081779423
$ws.Cells.Item(2, 12).Value = 1.065641108326689
$ws.Cells.Item(2, 13).Value = 1.071606285526605
$ws.Cells.Item(2, 14).Value = 1.056483135711364

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.051005769824119
$ws.Cells.Item(3, 4).Value = 1.055686556673284
$ws.Cells.Item(3, 5).Value = 1.063917271918247
$ws.Cells.Item(3, 6).Value = 1.069872227999325
$ws.Cells.Item(3, 9).Value = 1.044758878757909
$ws.Cells.Item(3, 10).Value = 1.055690472685513
$ws.Cells.Item(3, 11).Value = 1.058241898211722
$ws.Cells.Item(3, 12).Value = 1.066451765450334
$ws.Cells.Item(3, 13).Value = 1.072391853257151
$ws.Cells.Item(3, 14).Value = 1.057189674246116

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.051689741089316
$ws.Cells.Item(4, 4).Value = 1.056209649899412
$ws.Cells.Item(4, 5).Value = 1.064562430404167
$ws.Cells.Item(4, 6).Value = 1.070500159464865
$ws.Cells.Item(4, 9).Value = 1.044910147673852
$ws.Cells.Item(4, 10).Value = 1.056146703375557
$ws.Cells.Item(4, 11).Value = 1.058643558337905
$ws.Cells.Item(4, 12).Value = 1.066976246688129
$ws.Cells.Item(4, 13).Value = 1.072899897535925
$ws.Cells.Item(4, 14).Value = 1.057646552836071

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.051977353462623
$ws.Cells.Item(5, 4).Value = 1.056429547047916
$ws.Cells.Item(5, 5).Value = 1.064833785712013
$ws.Cells.Item(5, 6).Value = 1.070764225485134
$ws.Cells.Item(5, 9).Value = 1.044973350984675
$ws.Cells.Item(5, 10).Value = 1.05633843073768
$ws.Cells.Item(5, 11).Value = 1.058812248830436
$ws.Cells.Item(5, 12).Value = 1.067196721668583
$ws.Cells.Item(5, 13).Value = 1.073113413884946
$ws.Cells.Item(5, 14).Value = 1.057838552473039

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.052025648994053
$ws.Cells.Item(6, 4).Value = 1.056466468061576
$ws.Cells.Item(6, 5).Value = 1.064879355104284
$ws.Cells.Item(6, 6).Value = 1.070808568212461
$ws.Cells.Item(6, 9).Value = 1.044983940207922
$ws.Cells.Item(6, 10).Value = 1.056370618367027
$ws.Cells.Item(6, 11).Value = 1.058840562870006
$ws.Cells.Item(6, 12).Value = 1.067233739376953
$ws.Cells.Item(6, 13).Value = 1.073149260360428
$ws.Cells.Item(6, 14).Value = 1.057870785812511

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.051693583901796
$ws.Cells.Item(7, 4).Value = 1.056212588221174
$ws.Cells.Item(7, 5).Value = 1.0645660557545
$ws.Cells.Item(7, 6).Value = 1.070503687600359
$ws.Cells.Item(7, 9).Value = 1.044910993733123
$ws.Cells.Item(7, 10).Value = 1.056149265529598
$ws.Cells.Item(7, 11).Value = 1.058645813047192
$ws.Cells.Item(7, 12).Value = 1.066979192752923
$ws.Cells.Item(7, 13).Value = 1.072902750809938
$ws.Cells.Item(7, 14).Value = 1.057649118628665

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.050306208000286
$ws.Cells.Item(8, 4).Value = 1.055151318928828
$ws.Cells.Item(8, 5).Value = 1.063257626668787
$ws.Cells.Item(8, 6).Value = 1.069230048408936
$ws.Cells.Item(8, 9).Value = 1.044602779512546
$ws.Cells.Item(8, 10).Value = 1.055223437016438
$ws.Cells.Item(8, 11).Value = 1.057830371425432
$ws.Cells.Item(8, 12).Value = 1.065915087430216
$ws.Cells.Item(8, 13).Value = 1.071871827428743
$ws.Cells.Item(8, 14).Value = 1.056721975332821

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.047863757832967
$ws.Cells.Item(9, 4).Value = 1.053280934244155
$ws.Cells.Item(9, 5).Value = 1.060956242959215
$ws.Cells.Item(9, 6).Value = 1.066988470420456
$ws.Cells.Item(9, 9).Value = 1.044047218355241
$ws.Cells.Item(9, 10).Value = 1.053589744082132
$ws.Cells.Item(9, 11).Value = 1.056388151070225
$ws.Cells.Item(9, 12).Value = 1.064039503203817
$ws.Cells.Item(9, 13).Value = 1.070053166541676
$ws.Cells.Item(9, 14).Value = 1.055085962367162

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.046237002987379
$ws.Cells.Item(10, 4).Value = 1.052033864831091
$ws.Cells.Item(10, 5).Value = 1.059424875380455
$ws.Cells.Item(10, 6).Value = 1.065495971005864
$ws.Cells.Item(10, 9).Value = 1.043668499192024
$ws.Cells.Item(10, 10).Value = 1.052499121364212
$ws.Cells.Item(10, 11).Value = 1.055423137721644
$ws.Cells.Item(10, 12).Value = 1.062788815672735
$ws.Cells.Item(10, 13).Value = 1.068839391530578
$ws.Cells.Item(10, 14).Value = 1.053993790839886

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.045532967886565
$ws.Cells.Item(11, 4).Value = 1.051493845368318
$ws.Cells.Item(11, 5).Value = 1.05876246920732
$ws.Cells.Item(11, 6).Value = 1.064850158551016
$ws.Cells.Item(11, 9).Value = 1.043502535086889
$ws.Cells.Item(11, 10).Value = 1.052026520748512
$ws.Cells.Item(11, 11).Value = 1.05500444525814
$ws.Cells.Item(11, 12).Value = 1.062247190120597
$ws.Cells.Item(11, 13).Value = 1.068313504558313
$ws.Cells.Item(11, 14).Value = 1.053520519077104

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.045271512046621
$ws.Cells.Item(12, 4).Value = 1.051293254304076
$ws.Cells.Item(12, 5).Value = 1.058516525692559
$ws.Cells.Item(12, 6).Value = 1.064610343247402
$ws.Cells.Item(12, 9).Value = 1.043440592054523
$ws.Cells.Item(12, 10).Value = 1.051850923048979
$ws.Cells.Item(12, 11).Value = 1.054848799227733
$ws.Cells.Item(12, 12).Value = 1.062045996276631
$ws.Cells.Item(12, 13).Value = 1.068118120016629
$ws.Cells.Item(12, 14).Value = 1.053344672008696

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.045327592761922
$ws.Cells.Item(13, 4).Value = 1.051336281927963
$ws.Cells.Item(13, 5).Value = 1.058569276712926
$ws.Cells.Item(13, 6).Value = 1.064661781347736
$ws.Cells.Item(13, 9).Value = 1.043453892460577
$ws.Cells.Item(13, 10).Value = 1.051888591734566
$ws.Cells.Item(13, 11).Value = 1.054882191478609
$ws.Cells.Item(13, 12).Value = 1.062089153484217
$ws.Cells.Item(13, 13).Value = 1.068160032769123
$ws.Cells.Item(13, 14).Value = 1.053382394188137

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.045511354760832
$ws.Cells.Item(14, 4).Value = 1.051477264521433
$ws.Cells.Item(14, 5).Value = 1.058742137322413
$ws.Cells.Item(14, 6).Value = 1.064830333948295
$ws.Cells.Item(14, 9).Value = 1.043497420904635
$ws.Cells.Item(14, 10).Value = 1.052012006862656
$ws.Cells.Item(14, 11).Value = 1.054991582055995
$ws.Cells.Item(14, 12).Value = 1.062230559575999
$ws.Cells.Item(14, 13).Value = 1.06829735496087
$ws.Cells.Item(14, 14).Value = 1.053505984579864

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.045624583767442
$ws.Cells.Item(15, 4).Value = 1.051564128070041
$ws.Cells.Item(15, 5).Value = 1.058848656197779
$ws.Cells.Item(15, 6).Value = 1.064934193824412
$ws.Cells.Item(15, 9).Value = 1.0435242009265
$ws.Cells.Item(15, 10).Value = 1.052088040012183
$ws.Cells.Item(15, 11).Value = 1.055058964651015
$ws.Cells.Item(15, 12).Value = 1.062317683224998
$ws.Cells.Item(15, 13).Value = 1.068381957522287
$ws.Cells.Item(15, 14).Value = 1.053582125705185

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.046283735071972
$ws.Cells.Item(16, 4).Value = 1.052069703571967
$ws.Cells.Item(16, 5).Value = 1.059468851599082
$ws.Cells.Item(16, 6).Value = 1.065538840979942
$ws.Cells.Item(16, 9).Value = 1.043679472074747
$ws.Cells.Item(16, 10).Value = 1.052530478887862
$ws.Cells.Item(16, 11).Value = 1.05545090740421
$ws.Cells.Item(16, 12).Value = 1.062824760149622
$ws.Cells.Item(16, 13).Value = 1.068874286338332
$ws.Cells.Item(16, 14).Value = 1.054025192894816

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.046697299771226
$ws.Cells.Item(17, 4).Value = 1.052386830173121
$ws.Cells.Item(17, 5).Value = 1.059858067867217
$ws.Cells.Item(17, 6).Value = 1.065918241002655
$ws.Cells.Item(17, 9).Value = 1.043776340709354
$ws.Cells.Item(17, 10).Value = 1.052807914644985
$ws.Cells.Item(17, 11).Value = 1.055696539447
$ws.Cells.Item(17, 12).Value = 1.063142817798324
$ws.Cells.Item(17, 13).Value = 1.069183027401295
$ws.Cells.Item(17, 14).Value = 1.054303022642537

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.046938559845541
$ws.Cells.Item(18, 4).Value = 1.052571801792935
$ws.Cells.Item(18, 5).Value = 1.060085157292467
$ws.Cells.Item(18, 6).Value = 1.06613958195774
$ws.Cells.Item(18, 9).Value = 1.043832651847614
$ws.Cells.Item(18, 10).Value = 1.052969704069139
$ws.Cells.Item(18, 11).Value = 1.05583973183326
$ws.Cells.Item(18, 12).Value = 1.063328328697853
$ws.Cells.Item(18, 13).Value = 1.069363080471198
$ws.Cells.Item(18, 14).Value = 1.05446504182623

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.047020829196303
$ws.Cells.Item(19, 4).Value = 1.052634871810271
$ws.Cells.Item(19, 5).Value = 1.060162600123159
$ws.Cells.Item(19, 6).Value = 1.066215060835198
$ws.Cells.Item(19, 9).Value = 1.043851820134938
$ws.Cells.Item(19, 10).Value = 1.053024864276712
$ws.Cells.Item(19, 11).Value = 1.05588854304563
$ws.Cells.Item(19, 12).Value = 1.063391581962388
$ws.Cells.Item(19, 13).Value = 1.069424468747926
$ws.Cells.Item(19, 14).Value = 1.054520280367625

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.046652924595235
$ws.Cells.Item(20, 4).Value = 1.052352805790932
$ws.Cells.Item(20, 5).Value = 1.059816301797712
$ws.Cells.Item(20, 6).Value = 1.065877530483146
$ws.Cells.Item(20, 9).Value = 1.043765967343924
$ws.Cells.Item(20, 10).Value = 1.052778151944406
$ws.Cells.Item(20, 11).Value = 1.055670193789624
$ws.Cells.Item(20, 12).Value = 1.063108693917045
$ws.Cells.Item(20, 13).Value = 1.069149905544413
$ws.Cells.Item(20, 14).Value = 1.05427321767551

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.045457239937778
$ws.Cells.Item(21, 4).Value = 1.051435748760271
$ws.Cells.Item(21, 5).Value = 1.058691231316177
$ws.Cells.Item(21, 6).Value = 1.064780697520767
$ws.Cells.Item(21, 9).Value = 1.043484611047043
$ws.Cells.Item(21, 10).Value = 1.051975665637776
$ws.Cells.Item(21, 11).Value = 1.054959372704333
$ws.Cells.Item(21, 12).Value = 1.062188919278743
$ws.Cells.Item(21, 13).Value = 1.068256918283128
$ws.Cells.Item(21, 14).Value = 1.053469591746277

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.044705778790083
$ws.Cells.Item(22, 4).Value = 1.050859137036535
$ws.Cells.Item(22, 5).Value = 1.057984454704283
$ws.Cells.Item(22, 6).Value = 1.064091469987741
$ws.Cells.Item(22, 9).Value = 1.043305995451432
$ws.Cells.Item(22, 10).Value = 1.051470805408719
$ws.Cells.Item(22, 11).Value = 1.054511728263549
$ws.Cells.Item(22, 12).Value = 1.061610562838606
$ws.Cells.Item(22, 13).Value = 1.067695191788222
$ws.Cells.Item(22, 14).Value = 1.052964014557786

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.045104112874181
$ws.Cells.Item(23, 4).Value = 1.051164811618777
$ws.Cells.Item(23, 5).Value = 1.058359073315116
$ws.Cells.Item(23, 6).Value = 1.064456804898869
$ws.Cells.Item(23, 9).Value = 1.043400845479627
$ws.Cells.Item(23, 10).Value = 1.051738470180278
$ws.Cells.Item(23, 11).Value = 1.054749101421613
$ws.Cells.Item(23, 12).Value = 1.061917165839312
$ws.Cells.Item(23, 13).Value = 1.067992999044959
$ws.Cells.Item(23, 14).Value = 1.053232059444023

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.046672975719543
$ws.Cells.Item(24, 4).Value = 1.052368179953836
$ws.Cells.Item(24, 5).Value = 1.059835173882259
$ws.Cells.Item(24, 6).Value = 1.065895925680609
$ws.Cells.Item(24, 9).Value = 1.043770655210452
$ws.Cells.Item(24, 10).Value = 1.052791600533174
$ws.Cells.Item(24, 11).Value = 1.055682098507259
$ws.Cells.Item(24, 12).Value = 1.063124113050951
$ws.Cells.Item(24, 13).Value = 1.069164871979625
$ws.Cells.Item(24, 14).Value = 1.054286685362816

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.048494918120619
$ws.Cells.Item(25, 4).Value = 1.053764502980591
$ws.Cells.Item(25, 5).Value = 1.061550699701305
$ws.Cells.Item(25, 6).Value = 1.06756764299297
$ws.Cells.Item(25, 9).Value = 1.044192316890788
$ws.Cells.Item(25, 10).Value = 1.054012358300291
$ws.Cells.Item(25, 11).Value = 1.056761624941066
$ws.Cells.Item(25, 12).Value = 1.064524441739928
$ws.Cells.Item(25, 13).Value = 1.070523571993133
$ws.Cells.Item(25, 14).Value = 1.055509176745985
